$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") for rows 2..126 from 45190 to 45192
for ($r = 2; $r -le 126; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45190) {
        $cell.Value2 = 45192
    }
}

# Touch row 126's height so it gets serialized with an explicit custom height,
# matching Excel's behaviour when a row is edited.
$ws.Rows.Item(126).RowHeight = 15

# Append a new row (127) with new data
$row = 127
$ws.Cells.Item($row, 1).Value2 = "A 44580-2023"
$ws.Cells.Item($row, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item($row, 2).Value2 = 45189
$ws.Cells.Item($row, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item($row, 3).Value2 = 45192
$ws.Cells.Item($row, 4).Value2 = "SKÅNE LÄN"
$ws.Cells.Item($row, 5).Value2 = "SJÖBO"
$ws.Cells.Item($row, 6).Value2 = "Kommuner"
$ws.Cells.Item($row, 7).Value2 = 6.4
$ws.Cells.Item($row, 8).Value2 = 0
$ws.Cells.Item($row, 9).Value2 = 0
$ws.Cells.Item($row, 10).Value2 = 0
$ws.Cells.Item($row, 11).Value2 = 0
$ws.Cells.Item($row, 12).Value2 = 0
$ws.Cells.Item($row, 13).Value2 = 0
$ws.Cells.Item($row, 14).Value2 = 0
$ws.Cells.Item($row, 15).Value2 = 0
$ws.Cells.Item($row, 16).Value2 = 0
$ws.Cells.Item($row, 17).Value2 = 0

# R127 mirrors the wrap-text "empty" style used throughout column R
$rCell = $ws.Cells.Item($row, 18)
$rCell.WrapText = $true
$rCell.Value2 = "'"
